# "Bring in DeleteAsmtTest case."
# Replace the sample assignment-import test data with the rows pulled in
# from the DeleteAsmtTest case, then tidy up the view (autosize columns
# that now hold much longer strings, move the selection, set the page
# to portrait) the way a person doing this in the Excel UI would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "name_input"
$ws.Range("C1").Value = "due_date"
$ws.Range("D1").Value = "grade_type"
$ws.Range("E1").Value = "instructions"
$ws.Range("F1").Value = "post_sub_instructions"

# --- Row 2 ------------------------------------------------------------
# (A2 = 1 and C2 = 4/23/2018 are unchanged by this edit)
$ws.Range("B2").Value = "Ren_IP_0001"
$ws.Range("D2").Value = "Five Star"
$ws.Range("E2").Value = 'ins_symbol~!@#$%^&*()_+=-0987654321`{}:"|<>?][\'';,./'
$ws.Range("F2").Value = "This message is post submission instructions text."

# --- Row 3 ------------------------------------------------------------
# (A3 = 2 and C3 = 4/29/2018 are unchanged by this edit)
$ws.Range("B3").Value = "Ren_IP_0002"
$ws.Range("D3").Value = "AutoPass"
$ws.Range("E3").Value = "This message is instructions text."
$ws.Range("F3").Value = "PS_in Test long string less than 400 charactorsSteps to reproduce:1. Login bigben (https:bigben-moodle.youseeu.com) as educator-1.2. Select Course - bigbengenerallink.3. Create any type of project assignment, such as individual project and save it.4. Edit this assignment.5. Click the revord icon in Instructions.6. Click UPLOAD VIDEO button.7. Click SELECT FILE button to select file from local host."

# --- Resize columns now that E/F hold much longer text -----------------
$ws.Columns.Item(1).ColumnWidth = 8.996651785714286
$ws.Columns.Item(2).ColumnWidth = 17.711495535714285
$ws.Columns.Item(4).ColumnWidth = 18.570870535714285
$ws.Columns.Item(5).ColumnWidth = 64.14118303571429
$ws.Columns.Item(6).ColumnWidth = 45.141183035714285

# --- Page setup / selection, matching the saved view --------------------
$ws.PageSetup.Orientation = 1
$ws.Range("B4").Select()
